# Update "want to go" counts (F column) per commit diff
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F2").Value = 13796
$ws1.Range("F5").Value = 80
$ws1.Range("F6").Value = 806
$ws1.Range("F7").Value = 2201
$ws1.Range("F8").Value = 210
$ws1.Range("F10").Value = 120
$ws1.Range("F13").Value = 622
$ws1.Range("F14").Value = 471
$ws1.Range("F18").Value = 321
$ws1.Range("F19").Value = 910
$ws1.Range("F21").Value = 95
$ws1.Range("F25").Value = 120
$ws2.Range("F4").Value = 145
$ws2.Range("F6").Value = 141
$ws2.Range("F8").Value = 2157
$ws2.Range("F15").Value = 1922
$ws3.Range("F2").Value = 239
$ws3.Range("F3").Value = 211
$ws4.Range("F2").Value = 239
$ws4.Range("F3").Value = 13796
$ws4.Range("F6").Value = 80
$ws4.Range("F7").Value = 806
$ws4.Range("F10").Value = 2201
$ws4.Range("F11").Value = 211
$ws4.Range("F12").Value = 210
$ws4.Range("F14").Value = 120
$ws4.Range("F16").Value = 145
$ws4.Range("F19").Value = 141
$ws4.Range("F21").Value = 622
$ws4.Range("F22").Value = 471
$ws4.Range("F26").Value = 321
$ws4.Range("F27").Value = 910
$ws4.Range("F29").Value = 2157
$ws4.Range("F35").Value = 95
$ws4.Range("F41").Value = 120
$ws4.Range("F43").Value = 1922
